$wb = $excel.ActiveWorkbook

# --- Login sheet: add a third row of test credentials ---
$loginSheet = $wb.Worksheets.Item("Login")

$loginSheet.Range("A3").Value = "qa.test3@gmail.com"
$loginSheet.Range("B3").Value = "validPassword@456"

# New row styling: A3 -> 14pt, B3 -> 12pt, both in dark blue (FF333399)
$loginSheet.Range("A3").Font.Size = 14
$loginSheet.Range("A3").Font.Color = 10040115
$loginSheet.Range("B3").Font.Size = 12
$loginSheet.Range("B3").Font.Color = 10040115

$loginSheet.Rows.Item(3).RowHeight = 19

# --- Switch the active sheet/tab from contactus back to Login ---
$loginSheet.Activate()
$loginSheet.Range("A8").Select() | Out-Null
